$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 9836.5
$ws.Range("I62").Value = 8645.691999999999
$ws.Range("K62").Value = 8645.691999999999
$ws.Range("M62").Value = -8021.691999999999
$ws.Range("H65").Value = 9836.5
$ws.Range("I65").Value = 8645.691999999999
$ws.Range("K65").Value = 43228.45999999999
$ws.Range("M65").Value = -40108.45999999999
$ws.Range("H98").Value = 2438.8684
$ws.Range("J98").Value = 4668.1665
$ws.Range("L98").Value = 4668.1665
$ws.Range("N98").Value = -7664.1665
$ws.Range("H112").Value = 4522.5835
$ws.Range("J112").Value = 4608.943
$ws.Range("L112").Value = 13826.829
$ws.Range("N112").Value = -16042.829
$ws.Range("H113").Value = 3927.3333
$ws.Range("I113").Value = 3883
$ws.Range("J113").Value = 3949.5
$ws.Range("K113").Value = 3883
$ws.Range("L113").Value = 3949.5
$ws.Range("M113").Value = -629
$ws.Range("N113").Value = -10457.5
$ws.Range("H116").Value = 24479.75
$ws.Range("I116").Value = 26112
$ws.Range("J116").Value = 22381.143
$ws.Range("K116").Value = 26112
$ws.Range("L116").Value = 22381.143
$ws.Range("M116").Value = -22670
$ws.Range("N116").Value = -29265.143
$ws.Range("H122").Value = 2438.8684
$ws.Range("J122").Value = 4668.1665
$ws.Range("L122").Value = 14004.4995
$ws.Range("N122").Value = -18904.4995
$ws.Range("H125").Value = 4166.6665
$ws.Range("J125").Value = 4166.6665
$ws.Range("L125").Value = 37499.9985
$ws.Range("N125").Value = -42419.9985
$ws.Range("H138").Value = 1974.5641
$ws.Range("J138").Value = 3039.1428
$ws.Range("L138").Value = 9117.428400000001
$ws.Range("N138").Value = -19397.4284
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5355.778
$ws.Range("I45").Value = 2300
$ws.Range("J45").Value = 7800.4
$ws.Range("K45").Value = 2300
$ws.Range("L45").Value = 7800.4
$ws.Range("M45").Value = -1923
$ws.Range("N45").Value = -8554.4
$ws.Range("H61").Value = 7533.4707
$ws.Range("I61").Value = 1215.25
$ws.Range("K61").Value = 1215.25
$ws.Range("M61").Value = -1003.25
$ws.Range("H102").Value = 6203.56
$ws.Range("I102").Value = 6337.975
$ws.Range("K102").Value = 6337.975
$ws.Range("M102").Value = -4715.975
$ws.Range("H132").Value = 1157.875
$ws.Range("I132").Value = 929.1607
$ws.Range("K132").Value = 2787.4821
$ws.Range("M132").Value = -257.4821000000002
$ws.Range("H136").Value = 7533.4707
$ws.Range("I136").Value = 1215.25
$ws.Range("K136").Value = 3645.75
$ws.Range("M136").Value = -1095.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1780.0416
$ws.Range("I86").Value = 1729.3125
$ws.Range("J86").Value = 1881.5
$ws.Range("K86").Value = 1729.3125
$ws.Range("L86").Value = 1881.5
$ws.Range("M86").Value = -606.3125
$ws.Range("N86").Value = -4127.5
$ws.Range("H89").Value = 1780.0416
$ws.Range("I89").Value = 1729.3125
$ws.Range("J89").Value = 1881.5
$ws.Range("K89").Value = 8646.5625
$ws.Range("L89").Value = 9407.5
$ws.Range("M89").Value = -3030.5625
$ws.Range("N89").Value = -20639.5
$ws.Range("H94").Value = 6560.3184
$ws.Range("I94").Value = 8665.857
$ws.Range("J94").Value = 2875.625
$ws.Range("K94").Value = 8665.857
$ws.Range("L94").Value = 2875.625
$ws.Range("M94").Value = -8214.857
$ws.Range("N94").Value = -3777.625
$ws.Range("H123").Value = 28000
$ws.Range("I123").Value = 10000
$ws.Range("J123").Value = 100000
$ws.Range("K123").Value = 10000
$ws.Range("L123").Value = 100000
$ws.Range("M123").Value = -5100
$ws.Range("N123").Value = -109800
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1306.4286
$ws.Range("I22").Value = 281.33334
$ws.Range("K22").Value = 281.33334
$ws.Range("M22").Value = 68.66665999999998
$ws.Range("H105").Value = 1640.3572
$ws.Range("I105").Value = 1634.7273
$ws.Range("J105").Value = 1661
$ws.Range("K105").Value = 1634.7273
$ws.Range("L105").Value = 1661
$ws.Range("M105").Value = 112.2727
$ws.Range("N105").Value = -5155
$ws.Range("H132").Value = 67679.336
$ws.Range("I132").Value = 67679.336
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 203038.008
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -200508.008
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 2248.543
$ws.Range("I134").Value = 2204.16
$ws.Range("K134").Value = 6612.48
$ws.Range("M134").Value = -4077.48
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 26999.8
$ws.Range("J32").Value = 31249.75
$ws.Range("L32").Value = 93749.25
$ws.Range("N32").Value = -94315.25
$ws.Range("H134").Value = 1399.6842
$ws.Range("J134").Value = 3500
$ws.Range("L134").Value = 10500
$ws.Range("N134").Value = -20640
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7611.0454
$ws.Range("I80").Value = 5946.3335
$ws.Range("K80").Value = 5946.3335
$ws.Range("M80").Value = -4948.3335
$ws.Range("H83").Value = 7611.0454
$ws.Range("I83").Value = 5946.3335
$ws.Range("K83").Value = 29731.6675
$ws.Range("M83").Value = -24739.6675
$ws.Range("H97").Value = 1661.0555
$ws.Range("I97").Value = 1713
$ws.Range("J97").Value = 1619.5
$ws.Range("K97").Value = 1713
$ws.Range("L97").Value = 1619.5
$ws.Range("M97").Value = -1217
$ws.Range("N97").Value = -2611.5
$ws.Range("H102").Value = 22863.25
$ws.Range("J102").Value = 1220.2222
$ws.Range("L102").Value = 1220.2222
$ws.Range("N102").Value = -4464.2222
$ws.Range("H113").Value = 2171.4666
$ws.Range("I113").Value = 1990.1428
$ws.Range("J113").Value = 2330.125
$ws.Range("K113").Value = 1990.1428
$ws.Range("L113").Value = 2330.125
$ws.Range("M113").Value = 179.8571999999999
$ws.Range("N113").Value = -6670.125
$ws.Range("H132").Value = 3592.15
$ws.Range("I132").Value = 3588.2222
$ws.Range("J132").Value = 3627.5
$ws.Range("K132").Value = 10764.6666
$ws.Range("L132").Value = 10882.5
$ws.Range("M132").Value = -8234.6666
$ws.Range("N132").Value = -15942.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3058.5652
$ws.Range("J93").Value = 4737.5
$ws.Range("L93").Value = 4737.5
$ws.Range("N93").Value = -7233.5
$ws.Range("H122").Value = 5300
$ws.Range("I122").Value = 5666.6665
$ws.Range("K122").Value = 16999.9995
$ws.Range("M122").Value = -14549.9995
$ws.Range("H136").Value = 6605
$ws.Range("I136").Value = 6466.7
$ws.Range("K136").Value = 19400.1
$ws.Range("M136").Value = -16850.1
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 14716.386
$ws.Range("I132").Value = 17817.25
$ws.Range("J132").Value = 4221.154
$ws.Range("K132").Value = 53451.75
$ws.Range("L132").Value = 12663.462
$ws.Range("M132").Value = -50921.75
$ws.Range("N132").Value = -17723.462
